# Replace the 15 lattice-multiplication problems (3 columns x 5 rows) in the
# single table of the document with a new set of problems, keeping the
# existing table/cell/run formatting (the sz=32 run properties, cell
# widths, etc.) untouched and only rewriting each cell's text runs.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New "A x B" problems, listed in row-major cell order (row 1 col 1..3,
# row 2 col 1..3, ...). Each cell's five printed lines are derived from A
# (the two-digit multiplicand placed down the left side of the lattice)
# and B (the two-digit multiplier placed across the top):
#   line1 = "A x B"
#   line2 = "  b1    b2"     (the two digits of B, spaced out)
#   line3 = "  ----"
#   line4 = "a1|    |"       (first digit of A)
#   line5 = "a2|    |"       (second digit of A)
$problems = @(
    @("11", "62"), @("20", "57"), @("76", "84"),
    @("47", "33"), @("85", "75"), @("62", "58"),
    @("67", "60"), @("92", "20"), @("50", "53"),
    @("85", "63"), @("27", "21"), @("56", "35"),
    @("91", "51"), @("33", "89"), @("44", "61")
)

$rows = 5
$cols = 3
$index = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $a = $problems[$index][0]
        $b = $problems[$index][1]
        $index++

        $a1 = $a.Substring(0, 1)
        $a2 = $a.Substring(1, 1)
        $b1 = $b.Substring(0, 1)
        $b2 = $b.Substring(1, 1)

        $line1 = "$a x $b"
        $line2 = "  $b1    $b2"
        $line3 = "  ----"
        $line4 = "$a1|    |"
        $line5 = "$a2|    |"

        $runInner = "<w:t>$line1</w:t><w:br/>" +
            "<w:t xml:space=`"preserve`">$line2</w:t><w:br/>" +
            "<w:t xml:space=`"preserve`">$line3</w:t><w:br/>" +
            "<w:t>$line4</w:t><w:br/>" +
            "<w:t>$line5</w:t>"

        $cellXml = "<?xml version=`"1.0`" standalone=`"yes`"?>" +
            "<?mso-application progid=`"Word.Document`"?>" +
            "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
            "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
            "<pkg:xmlData>" +
            "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
            "<w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>$runInner</w:r></w:p></w:body>" +
            "</w:document></pkg:xmlData></pkg:part></pkg:package>"

        $cell = $t.Cell($r, $c)
        # InsertXML replaces exactly the contents of the Range it is called
        # on, so scoping to the cell's own Range rewrites only that cell
        # (and lets us control the xml:space="preserve" attribute exactly,
        # which a plain Range.Text assignment does not).
        $cell.Range.InsertXML($cellXml)
    }
}
